$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AC001"
$ws.Range("B2").Value = "A001"
$ws.Range("C2").Value = "MR3c45380b"
$ws.Range("D2").Value = "Diagnosis done"
$ws.Range("E2").Value = "Pending"

$ws.Range("H10").Select()
